$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row for the research paper on Coconut Mill Effluent pretreatment:
# add the "Measured Parameters" value in column D of the new row 6.
$ws.Range("D6").Value = "Glycerol/Free fatty acids/Lipids/Proteins/Reducing sugars/Biogas/O&G/COD"

# Row 6 uses a taller row height (wrapped text), matching the authored row.
$ws.Rows.Item(6).RowHeight = 30

# Column D widened to comfortably fit the new, longer text.
$ws.Columns.Item(4).ColumnWidth = 48.42

# Update the active selection as left by the author after the edit.
$ws.Range("G6").Select() | Out-Null
